$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 5740.6816
$ws.Range("I19").Value = 2330.3
$ws.Range("J19").Value = 8582.666999999999
$ws.Range("K19").Value = 2330.3
$ws.Range("L19").Value = 8582.666999999999
$ws.Range("M19").Value = -2155.3
$ws.Range("N19").Value = -8932.666999999999

$ws.Range("H28").Value = 1522.1875
$ws.Range("I28").Value = 1254.5
$ws.Range("K28").Value = 1254.5
$ws.Range("M28").Value = -769.5

$ws.Range("H53").Value = 4297.8335
$ws.Range("I53").Value = 4841.5884
$ws.Range("J53").Value = 2977.2856
$ws.Range("K53").Value = 4841.5884
$ws.Range("L53").Value = 2977.2856
$ws.Range("M53").Value = -4204.5884
$ws.Range("N53").Value = -4251.2856

$ws.Range("H70").Value = 2596.1667
$ws.Range("J70").Value = 3984.8572
$ws.Range("L70").Value = 11954.5716
$ws.Range("N70").Value = -12494.5716

$ws.Range("H73").Value = 2596.1667
$ws.Range("J73").Value = 3984.8572
$ws.Range("L73").Value = 11954.5716
$ws.Range("N73").Value = -13826.5716

$ws.Range("H88").Value = 18220.16
$ws.Range("I88").Value = 1719.4
$ws.Range("J88").Value = 29220.666
$ws.Range("K88").Value = 1719.4
$ws.Range("L88").Value = 29220.666
$ws.Range("M88").Value = -1313.4
$ws.Range("N88").Value = -30032.666

$ws.Range("H91").Value = 18220.16
$ws.Range("I91").Value = 1719.4
$ws.Range("J91").Value = 29220.666
$ws.Range("K91").Value = 1719.4
$ws.Range("L91").Value = 29220.666
$ws.Range("M91").Value = -315.4000000000001
$ws.Range("N91").Value = -32028.666

$ws.Range("H98").Value = 2600.3
$ws.Range("I98").Value = 2600.3
$ws.Range("K98").Value = 2600.3
$ws.Range("M98").Value = -1102.3

$ws.Range("H107").Value = 540.15
$ws.Range("I107").Value = 459.69232
$ws.Range("K107").Value = 459.69232
$ws.Range("M107").Value = 1460.30768

$ws.Range("H113").Value = 5494.9
$ws.Range("I113").Value = 5494.9
$ws.Range("K113").Value = 5494.9
$ws.Range("M113").Value = -2240.9

$ws.Range("H116").Value = 6281.5356
$ws.Range("I116").Value = 5905.1763
$ws.Range("K116").Value = 5905.1763
$ws.Range("M116").Value = -2463.1763

$ws.Range("H122").Value = 2600.3
$ws.Range("I122").Value = 2600.3
$ws.Range("K122").Value = 7800.900000000001
$ws.Range("M122").Value = -5350.900000000001

$ws.Range("H138").Value = 4489.5566
$ws.Range("J138").Value = 5155.897
$ws.Range("L138").Value = 15467.691
$ws.Range("N138").Value = -25747.691

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2530.7222
$ws.Range("I2").Value = 2409.5625
$ws.Range("K2").Value = 2409.5625
$ws.Range("M2").Value = -2296.5625

$ws.Range("H32").Value = 13182.23
$ws.Range("I32").Value = 8484.508
$ws.Range("J32").Value = 21906.572
$ws.Range("K32").Value = 8484.508
$ws.Range("L32").Value = 21906.572
$ws.Range("M32").Value = -8197.508
$ws.Range("N32").Value = -22480.572

$ws.Range("H61").Value = 3670.1807
$ws.Range("I61").Value = 2898.9246
$ws.Range("K61").Value = 2898.9246
$ws.Range("M61").Value = -2686.9246

$ws.Range("H102").Value = 2206.9697
$ws.Range("I102").Value = 2088.5312
$ws.Range("K102").Value = 2088.5312
$ws.Range("M102").Value = -466.5311999999999

$ws.Range("H110").Value = 33814.535
$ws.Range("I110").Value = 37421.56
$ws.Range("K110").Value = 37421.56
$ws.Range("M110").Value = -35376.56

$ws.Range("H116").Value = 2530.7222
$ws.Range("I116").Value = 2409.5625
$ws.Range("K116").Value = 2409.5625
$ws.Range("M116").Value = -115.5625

$ws.Range("H136").Value = 3670.1807
$ws.Range("I136").Value = 2898.9246
$ws.Range("K136").Value = 8696.773799999999
$ws.Range("M136").Value = -6146.773799999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2530.7222
$ws.Range("I3").Value = 2409.5625
$ws.Range("K3").Value = 2409.5625
$ws.Range("M3").Value = -2295.5625

$ws.Range("H86").Value = 4324
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 4324
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()

$ws.Range("H94").Value = 2386.5334
$ws.Range("I94").Value = 2184.1
$ws.Range("J94").Value = 2791.4
$ws.Range("K94").Value = 2184.1
$ws.Range("L94").Value = 2791.4
$ws.Range("M94").Value = -1733.1
$ws.Range("N94").Value = -3693.4

$ws.Range("H107").Value = 113749.5
$ws.Range("I107").Value = 113749.5
$ws.Range("K107").Value = 113749.5
$ws.Range("M107").Value = -111829.5

$ws.Range("H134").Value = 3010.375
$ws.Range("I134").Value = 1807.1936
$ws.Range("K134").Value = 5421.5808
$ws.Range("M134").Value = -2886.5808

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2710.606
$ws.Range("I58").Value = 1504
$ws.Range("J58").Value = 3992.625
$ws.Range("K58").Value = 1504
$ws.Range("L58").Value = 3992.625
$ws.Range("M58").Value = -1301
$ws.Range("N58").Value = -4398.625

$ws.Range("H99").Value = 4786.3125
$ws.Range("I99").Value = 4049.1943
$ws.Range("K99").Value = 4049.1943
$ws.Range("M99").Value = -2551.1943

$ws.Range("H126").Value = 4786.3125
$ws.Range("I126").Value = 4049.1943
$ws.Range("K126").Value = 12147.5829
$ws.Range("M126").Value = -9677.582900000001

$ws.Range("H132").Value = 7080.8247
$ws.Range("I132").Value = 4971.9185
$ws.Range("J132").Value = 19997.875
$ws.Range("K132").Value = 14915.7555
$ws.Range("L132").Value = 59993.625
$ws.Range("M132").Value = -12385.7555
$ws.Range("N132").Value = -65053.625

$ws.Range("H136").Value = 2710.606
$ws.Range("I136").Value = 1504
$ws.Range("J136").Value = 3992.625
$ws.Range("K136").Value = 4512
$ws.Range("L136").Value = 11977.875
$ws.Range("M136").Value = -1962
$ws.Range("N136").Value = -17077.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 1575.1818
$ws.Range("I14").Value = 1575.1818
$ws.Range("K14").Value = 4725.5454
$ws.Range("M14").Value = -4552.5454

$ws.Range("H113").Value = 77433.664
$ws.Range("I113").Value = 4480
$ws.Range("J113").Value = 129543.43
$ws.Range("K113").Value = 13440
$ws.Range("L113").Value = 388630.29
$ws.Range("M113").Value = -11270
$ws.Range("N113").Value = -392970.29

$ws.Range("H118").Value = 1999
$ws.Range("I118").Value = 1999
$ws.Range("K118").Value = 5997
$ws.Range("M118").Value = -4754

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 33032.5
$ws.Range("J52").Value = 33033
$ws.Range("L52").Value = 33033
$ws.Range("N52").Value = -33551

$ws.Range("H55").Value = 8107.5
$ws.Range("I55").Value = 2476.6667
$ws.Range("J55").Value = 25000
$ws.Range("K55").Value = 2476.6667
$ws.Range("L55").Value = 25000
$ws.Range("M55").Value = -2149.6667
$ws.Range("N55").Value = -25654

$ws.Range("H70").Value = 4690
$ws.Range("I70").Value = 4491.4375
$ws.Range("K70").Value = 4491.4375
$ws.Range("M70").Value = -4221.4375

$ws.Range("H73").Value = 4690
$ws.Range("I73").Value = 4491.4375
$ws.Range("K73").Value = 4491.4375
$ws.Range("M73").Value = -3555.4375

$ws.Range("H80").Value = 18199.111
$ws.Range("I80").Value = 9723.5
$ws.Range("K80").Value = 9723.5
$ws.Range("M80").Value = -8725.5

$ws.Range("H83").Value = 18199.111
$ws.Range("I83").Value = 9723.5
$ws.Range("K83").Value = 48617.5
$ws.Range("M83").Value = -43625.5

$ws.Range("H97").Value = 880.7727
$ws.Range("I97").Value = 693.2222
$ws.Range("K97").Value = 693.2222
$ws.Range("M97").Value = -197.2222

$ws.Range("H107").Value = 796
$ws.Range("J107").Value = 795
$ws.Range("L107").Value = 795
$ws.Range("N107").Value = -4635

$ws.Range("H109").Value = 25980.666
$ws.Range("J109").Value = 25980.666
$ws.Range("L109").Value = 25980.666
$ws.Range("N109").Value = -28060.666

$ws.Range("H113").Value = 16293.875
$ws.Range("I113").Value = 16293.875
$ws.Range("K113").Value = 16293.875
$ws.Range("M113").Value = -14123.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 15874023
$ws.Range("I81").Value = 15874023
$ws.Range("K81").Value = 31748046
$ws.Range("M81").Value = -31746985

$ws.Range("H84").Value = 15874023
$ws.Range("I84").Value = 15874023
$ws.Range("K84").Value = 158740230
$ws.Range("M84").Value = -158734926

$ws.Range("H132").Value = 6873.5884
$ws.Range("I132").Value = 6486.5
$ws.Range("K132").Value = 19459.5
$ws.Range("M132").Value = -16929.5

$ws.Range("H136").Value = 22226560
$ws.Range("I136").Value = 31255490
$ws.Range("K136").Value = 93766470
$ws.Range("M136").Value = -93763920
